$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 11605
$ws1.Range("F3").Value = 11160
$ws1.Range("F6").Value = 1010
$ws1.Range("F8").Value = 68
$ws1.Range("F11").Value = 10684
$ws1.Range("F12").Value = 4131
$ws1.Range("F13").Value = 13
$ws1.Range("F15").Value = 2461
$ws1.Range("F16").Value = 814
$ws1.Range("F19").Value = 434
$ws1.Range("F20").Value = 11118
$ws1.Range("F21").Value = 10885
$ws1.Range("F26").Value = 27

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 11605
$ws4.Range("F3").Value = 11160
$ws4.Range("F6").Value = 1010
$ws4.Range("F8").Value = 68
$ws4.Range("F11").Value = 10684
$ws4.Range("F12").Value = 4131
$ws4.Range("F13").Value = 13
$ws4.Range("F16").Value = 814
$ws4.Range("F19").Value = 434
$ws4.Range("F20").Value = 11118
$ws4.Range("F21").Value = 10885
$ws4.Range("F26").Value = 27
